$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-missing "Households" counts (column C) for the
# "None of the Above" appliance rows (Survey 6, rows 23-29).
$ws.Range("C23").Value = 30
$ws.Range("C24").Value = 18
$ws.Range("C25").Value = 3
$ws.Range("C26").Value = 47
$ws.Range("C27").Value = 3
$ws.Range("C28").Value = 9
$ws.Range("C29").Value = 15

# Reflect the updated selection seen after the edit.
$ws.Range("E28").Select()
